$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 163729
$ws.Range("C4").Value = 154720
$ws.Range("C7").Value = 5.5
$ws.Range("C8").Value = 64.87
